# Add a new "backup" column (R) to the worksheet, reset a handful of
# previously non-zero "detect_structure" (Q) values back to 0, flip
# isPivot (O263) for one row, and append three new monthly rows
# (266-268) of trailing stock data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cell R1 = "backup", formatted like the rest of row 1 ---
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value = "backup"

# --- 2. Default every existing data row (2-265) to backup = 0 ---
$ws.Range("R2:R265").Value = 0

# --- 3. Rows whose backup value keeps the (previous) detect_structure value ---
$ws.Range("R52").Value = 2
$ws.Range("R59").Value = 2
$ws.Range("R225").Value = 1
$ws.Range("R245").Value = 1

# --- 4. Rows whose detect_structure (Q) value is reset back to 0 ---
$ws.Range("Q27").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("Q43").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("Q48").Value = 0
$ws.Range("Q53").Value = 0

# --- 5. Row 263 gets flagged as a pivot ---
$ws.Range("O263").Value = 1

# --- 6. Append three new trailing monthly rows, copying number formats
#        from the last existing row (265) first so columns A/others match ---
$ws.Range("A265:Q265").Copy()
$ws.Range("A266:Q266").PasteSpecial(-4122)
$ws.Range("A267:Q267").PasteSpecial(-4122)
$ws.Range("A268:Q268").PasteSpecial(-4122)

function Set-Row($r, $a, $b, $c, $d, $e, $f, $g, $h, $i, $j, $k, $l, $m, $n, $o, $p, $q) {
    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
    $ws.Range("I$r").Value = $i
    $ws.Range("J$r").Value = $j
    $ws.Range("K$r").Value = $k
    $ws.Range("L$r").Value = $l
    $ws.Range("M$r").Value = $m
    $ws.Range("N$r").Value = $n
    $ws.Range("O$r").Value = $o
    $ws.Range("P$r").Value = $p
    $ws.Range("Q$r").Value = $q
}

Set-Row 266 45474 123.1500015258789 128.8000030517578 112.5199966430664 123.9499969482422 123.9499969482422 809974836 2024 7 1 0 0 0 27 0 0 0
Set-Row 267 45505 124.5             125.6999969482422 111.9000015258789 116.5699996948242 116.5699996948242 469107715 2024 8 1 0 0 0 31 0 0 2
Set-Row 268 45536 116.9400024414062 117.4899978637695 103.620002746582  109.2200012207031 109.2200012207031 571925935 2024 9 1 0 0 0 35 0 0 0
